# Update the "Förändrad" date column (C) for rows 2-34 from 2024-12-16 (45642)
# to 2024-12-17 (45643), i.e. advance the serial date value by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 34; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45642) {
        $cell.Value = 45643
    }
}
